# Projektstrukturplan: renumber IDs
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Tabelle1")
$ws2 = $wb.Worksheets.Item("Tabelle2")
$ws3 = $wb.Worksheets.Item("Tabelle3")

# --- New named range: _xlnm.Print_Area_0 (sheet-scoped on Tabelle1) ---
$null = $ws1.Names.Add("_xlnm.Print_Area_0", "=Tabelle1!`$A`$1:`$R`$13")

# --- Row 3: fill in the previously-empty ID cells ---
$ws1.Range("B3").Value = 1
$ws1.Range("E3").Value = 2
$ws1.Range("H3").Value = 3
$ws1.Range("K3").Value = 4
$ws1.Range("N3").Value = 5
$ws1.Range("Q3").Value = 6

# --- Row 5: renumber ---
$ws1.Range("B5").Value = 10
$ws1.Range("E5").Value = 20
$ws1.Range("H5").Value = 30
$ws1.Range("K5").Value = 40
$ws1.Range("N5").Value = 50
$ws1.Range("Q5").Value = 60

# --- Row 7: renumber ---
$ws1.Range("B7").Value = 11
$ws1.Range("E7").Value = 21
$ws1.Range("H7").Value = 31
$ws1.Range("K7").Value = 41
$ws1.Range("N7").Value = 51
$ws1.Range("Q7").Value = 61

# --- Row 9: renumber ---
$ws1.Range("B9").Value = 12
$ws1.Range("E9").Value = 22
$ws1.Range("H9").Value = 32
$ws1.Range("K9").Value = 42
$ws1.Range("Q9").Value = 62

# --- Row 11: renumber ---
$ws1.Range("E11").Value = 23
$ws1.Range("H11").Value = 33
$ws1.Range("K11").Value = 43

# --- Row 13: renumber ---
$ws1.Range("H13").Value = 34
$ws1.Range("K13").Value = 44

# --- View state: zoom to 65% on every sheet, update selections ---
$null = $ws2.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 65
$null = $ws2.Range("A1").Select()

$null = $ws3.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 65
$null = $ws3.Range("A1").Select()

$null = $ws1.Activate()
$win = $excel.ActiveWindow
$win.Zoom = 65
$null = $ws1.Range("Q11").Select()
